$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) format, used to reset style
# after forcing text storage on numeric-looking D-column values.
$defaultStyle = $ws.Range("D4").Style

$ws.Range('D2').Value = '42.097.00'
$ws.Range('E2').Value = '  -4.15%  '
$ws.Range('D3').Value = '2.241.64'
$ws.Range('E3').Value = '  -4.84%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '232.97'
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Value = '  -3.32%  '
$ws.Range('E6').Value = '  -6.08%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '70.14'
$ws.Range('D7').Style = $defaultStyle
$ws.Range('E7').Value = '  -4.73%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  -7.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0993'
$ws.Range('D10').Style = $defaultStyle
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '58.37'
$ws.Range('D11').Style = $defaultStyle
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '35.61'
$ws.Range('D12').Style = $defaultStyle
$ws.Range('E12').Value = '  +5.38%  '
$ws.Range('E13').Value = '  -2.99%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.81'
$ws.Range('D14').Style = $defaultStyle
$ws.Range('E14').Value = '  -7.33%  '
$ws.Range('D15').Value = '2.574.46'
$ws.Range('E15').Value = '  -4.87%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.07'
$ws.Range('D16').Style = $defaultStyle
$ws.Range('E16').Value = '  -8.48%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.864'
$ws.Range('D17').Style = $defaultStyle
$ws.Range('E17').Value = '  -5.34%  '
$ws.Range('D18').Value = '2.239.75'
$ws.Range('E18').Value = '  -4.84%  '
$ws.Range('D19').Value = '41.996.60'
$ws.Range('E19').Value = '  -4.16%  '
$ws.Range('D20').Value = '0.0₃0983'
$ws.Range('E20').Value = '  -4.62%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.26'
$ws.Range('D21').Style = $defaultStyle
$ws.Range('E21').Value = '  -7.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '73.40'
$ws.Range('D22').Style = $defaultStyle
$ws.Range('E22').Value = '  -5.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '237.39'
$ws.Range('D23').Style = $defaultStyle
$ws.Range('E23').Value = '  -7.79%  '
$ws.Range('E24').Value = '  +5.84%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.65'
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('E27').Value = '  -5.95%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.04'
$ws.Range('D28').Style = $defaultStyle
$ws.Range('E28').Value = '  -5.88%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.16'
$ws.Range('D29').Style = $defaultStyle
$ws.Range('E29').Value = '  -5.57%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.16'
$ws.Range('D30').Style = $defaultStyle
$ws.Range('E30').Value = '  -5.19%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.66'
$ws.Range('D31').Style = $defaultStyle
$ws.Range('E32').Value = '  -7.75%  '
$ws.Range('E33').Value = '  -7.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.44'
$ws.Range('D34').Style = $defaultStyle
$ws.Range('E34').Value = '  -1.15%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0719'
$ws.Range('D35').Style = $defaultStyle
$ws.Range('E35').Value = '  -5.24%  '
$ws.Range('E36').Value = '  -8.31%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.61'
$ws.Range('D37').Style = $defaultStyle
$ws.Range('E37').Value = '  -5.47%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '22.29'
$ws.Range('D38').Style = $defaultStyle
$ws.Range('E38').Value = '  +16.62%  '
$ws.Range('E39').Value = '  -6.32%  '
$ws.Range('E40').Value = '  -5.98%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0267'
$ws.Range('D41').Style = $defaultStyle
$ws.Range('E41').Value = '  -3.91%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '67.46'
$ws.Range('D42').Style = $defaultStyle
$ws.Range('E42').Value = '  -1.20%  '
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.96'
$ws.Range('D43').Style = $defaultStyle
$ws.Range('E43').Value = '  -3.27%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '9.10'
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').Value = '  -2.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.101'
$ws.Range('D45').Style = $defaultStyle
$ws.Range('E45').Value = '  -9.84%  '
$ws.Range('E46').Value = '  -6.96%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.37'
$ws.Range('D48').Style = $defaultStyle
$ws.Range('E48').Value = '  -5.56%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.36'
$ws.Range('D49').Style = $defaultStyle
$ws.Range('E49').Value = '  +4.89%  '
$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '10.09'
$ws.Range('D50').Style = $defaultStyle
$ws.Range('E50').Value = '  +5.06%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.17'
$ws.Range('D51').Style = $defaultStyle
$ws.Range('E51').Value = '  -6.95%  '
